# Textbox response formatting fix
# Rename sheets and update task-order / stimulus filename cells.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (also updates sheet name / r:id mapping in workbook.xml) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16511688017936182"
$wb.Worksheets.Item(2).Name = "NB_TO-16511688044212458"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168804423249"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511688044790208"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511688045670817"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168801762368.csv"
$ws1.Range("B3").Value = "GNG_stims-16511688017779617.csv"
$ws1.Range("B4").Value = "go_stims-16511688017779617.csv"
$ws1.Range("B5").Value = "GNG_stims-16511688017936182.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16511688034182434.csv"
$ws2.Range("B3").Value = "ZB-match_3-16511688020225523.csv"
$ws2.Range("B4").Value = "TB-16511688044022477.csv"
$ws2.Range("B5").Value = "OB-16511688025447352.csv"
$ws2.Range("B6").Value = "ZB-match_4-16511688018899717.csv"
$ws2.Range("B7").Value = "ZB-match_4-16511688020069253.csv"
$ws2.Range("B8").Value = "OB-16511688022810228.csv"
$ws2.Range("B9").Value = "TB-16511688040432467.csv"
$ws2.Range("B10").Value = "TB-16511688037862477.csv"

# --- Sheet 3 (RS) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511688044460244.csv"
$ws4.Range("B3").Value = "ZM_stims-16511688044262452.csv"
$ws4.Range("B4").Value = "MM_stims-16511688044620302.csv"
$ws4.Range("B5").Value = "ZM_stims-16511688044470263.csv"
$ws4.Range("B6").Value = "MM_stims-16511688044780216.csv"
$ws4.Range("B7").Value = "ZM_stims-1651168804463024.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1651168804551253.csv"
$ws5.Range("B3").Value = "SAT_stims-16511688045100248.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168804525414.csv"
$ws5.Range("B5").Value = "SAT_stims-16511688044860225.csv"
